$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report was inserted above the current row 5, pushing
# the existing rows 5 and 6 down to 6 and 7 (dimension grows from T6 to T7).
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly entry.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44489
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101001
$ws.Range("J5").Value = "Arándano (blue)"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 9500
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9750
$ws.Range("Q5").Value = "`$/bandeja 2 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 4875
$ws.Range("T5").Value = 2
